$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source roster stores every cell (including NPIs, dates, zip-like
# PPG ids, etc.) as literal text, not numbers/dates. Force the whole
# A1:Q4 block to Text format BEFORE writing values so Excel doesn't
# auto-convert things like "08/01/2025" into a date serial or strip
# leading zeros from "041104".
$ws.Range("A1:Q4").NumberFormat = "@"

# ---- Header row (row 1) ----
$ws.Range("A1").Value = "Transaction Type"
$ws.Range("B1").Value = "Transaction Attribute"
$ws.Range("C1").Value = "Effective Date"
$ws.Range("D1").Value = "Termination Date"
$ws.Range("E1").Value = "Termination Reason"
$ws.Range("F1").Value = "Provider Name"
$ws.Range("G1").Value = "Provider NPI"
$ws.Range("H1").Value = "Provider Specialty"
$ws.Range("I1").Value = "State License"
$ws.Range("J1").Value = "Organization Name"
$ws.Range("K1").Value = "TIN"
$ws.Range("L1").Value = "Group NPI"
$ws.Range("M1").Value = "Complete Address"
$ws.Range("N1").Value = "Phone Number"
$ws.Range("O1").Value = "Fax Number"
$ws.Range("P1").Value = "PPG ID"
$ws.Range("Q1").Value = "Lines of Business(Medicare/Commercial/Medical)"

# ---- Row 2 (Cole) ----
$ws.Range("A2").Value = "Term"
$ws.Range("B2").Value = "Provider Name"
$ws.Range("C2").Value = "Information not found"
$ws.Range("D2").Value = "08/01/2025"
$ws.Range("E2").Value = "Voluntary"
$ws.Range("F2").Value = "Cole"
$ws.Range("G2").Value = "1222222250"
$ws.Range("H2").Value = "Pediatric Emergency Medicine"
$ws.Range("I2").Value = "Information not found"
$ws.Range("J2").Value = "Rchn"
$ws.Range("K2").Value = "821111113"
$ws.Range("L2").Value = "Information not found"
$ws.Range("M2").Value = "Information not found"
$ws.Range("N2").Value = "Information not found"
$ws.Range("O2").Value = "Information not found"
$ws.Range("P2").Value = "Information not found"
$ws.Range("Q2").Value = "FFS/PPO/ACO/HMO/Medi-Cal"

# ---- Row 3 (Cyrus) ----
$ws.Range("A3").Value = "Term"
$ws.Range("B3").Value = "Provider"
$ws.Range("C3").Value = "09/01/2025"
$ws.Range("D3").Value = "Information not found"
$ws.Range("E3").Value = "Information not found"
$ws.Range("F3").Value = "Cyrus"
$ws.Range("G3").Value = "1164444443"
$ws.Range("H3").Value = "Internal Medicine"
$ws.Range("I3").Value = "D66661"
$ws.Range("J3").Value = "Mercian"
$ws.Range("K3").Value = "458888885"
$ws.Range("L3").Value = "1999999997"
$ws.Range("M3").Value = "Information not found"
$ws.Range("N3").Value = "Information not found"
$ws.Range("O3").Value = "Information not found"
$ws.Range("P3").Value = "041104"
$ws.Range("Q3").Value = "Medicare, PPG#’s, Commercial HMO"

# ---- Row 4 (Paul) ----
$ws.Range("A4").Value = "Update"
$ws.Range("B4").Value = "Provider"
$ws.Range("C4").Value = "09/22/2025"
$ws.Range("D4").Value = "Information not found"
$ws.Range("E4").Value = "Information not found"
$ws.Range("F4").Value = "Paul"
$ws.Range("G4").Value = "Information not found"
$ws.Range("H4").Value = "Information not found"
$ws.Range("I4").Value = "Information not found"
$ws.Range("J4").Value = "Hilabs"
$ws.Range("K4").Value = "Information not found"
$ws.Range("L4").Value = "Information not found"
$ws.Range("M4").Value = "Information not found"
$ws.Range("N4").Value = "Information not found"
$ws.Range("O4").Value = "Information not found"
$ws.Range("P4").Value = "Information not found"
$ws.Range("Q4").Value = "Information not found"

# ---- Drop the now-unused Fax/Email columns (R:S) so the sheet's
#      dimension shrinks from A1:S4 down to A1:Q4 ----
$ws.Range("R1:S4").Delete()
